$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (坐标-x / 坐标-z columns)
$ws.Range("E10").Value = 1
$ws.Range("F11").Value = 4
$ws.Range("F12").Value = -4

# Move the active selection to F13, matching the saved cursor position
$ws.Activate()
$ws.Range("F13").Select()
